$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- First touch each brand-new text value in the exact order they were
# --- originally authored, so the shared-strings table comes out in the
# --- same order as the source workbook.
$ws.Range("D28").Value = "Logger implementiert."
$ws.Range("D24").Value = "Weekly-Summup-02 Meeting."
$ws.Range("D25").Value = "Zusammenfassung Weekly-Summup-02 erstellt."
$ws.Range("D26").Value = "Weekly-Summup-03 Meeting"
$ws.Range("D27").Value = "Zusammenfassung Weekly-Summup-03 erstellt."
$ws.Range("D30").Value = "Sprint Review #1"
$ws.Range("A32").Value = "Gesamt:"
$ws.Range("D29").Value = "Meeting mit Gruppe WIFSurvivors bzgl. Standardisierung der Szenen-JSON Datei."

# --- Now fill in the remaining cells for the new booking rows (24-30) ---
# (date cells get the raw serial number so no ad-hoc number format gets
# auto-created; the correct date display format is applied later via
# PasteSpecial of the formatting from row 23)
$ws.Range("A24").Value = 45593
$ws.Range("B24").Value = 1.5
$ws.Range("C24").Value = "Online-Meeting"

$ws.Range("A25").Value = 45593
$ws.Range("B25").Value = 0.5
$ws.Range("C25").Value = "Planung"

$ws.Range("A26").Value = 45600
$ws.Range("B26").Value = 1.5
$ws.Range("C26").Value = "Online-Meeting"

$ws.Range("A27").Value = 45600
$ws.Range("B27").Value = 0.5
$ws.Range("C27").Value = "Planung"

$ws.Range("A28").Value = 45600
$ws.Range("B28").Value = 8
$ws.Range("C28").Value = "Coding"

$ws.Range("A29").Value = 45601
$ws.Range("B29").Value = 1.5
$ws.Range("C29").Value = "Online-Meeting"

$ws.Range("A30").Value = 45603
$ws.Range("B30").Value = 0.5
$ws.Range("C30").Value = "Online-Meeting"

# --- Totals row (row 31 stays empty) ---
$ws.Range("B32").Formula = "=SUM(B7:B30)"

# --- Match formatting of the copied rows (date format for A, plain for B) ---
$ws.Range("A23:B23").Copy()
$ws.Range("A24:B30").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B1").Copy()
$ws.Range("A32:B32").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- Update selection / scroll position to match final saved view ---
$ws.Range("D30").Select()
$excel.ActiveWindow.ScrollRow = 7
